# Suite 1.xlsx - "driven framework" keyword-table expansion.
#
# The original sheet held a single data row (row 3): 1 | FileSystem |
# checkFileContainsKeyword | "" | "".  The edit blanks that row out and
# re-homes the FileSystem / checkFileContainsKeyword keyword call onto row
# 4, padding rows 3, 5 and 6 with empty cells so the used range grows from
# A1:E3 to A1:E6 (a 4-row driver table instead of a 1-row one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear out the old single data row and lay out the new blank rows ---
$ws.Range("A3:E3").Value = ""
$ws.Range("A4").Value = ""
$ws.Range("B4").Value = "FileSystem"
$ws.Range("C4").Value = "checkFileContainsKeyword"
$ws.Range("D4:E4").Value = ""
$ws.Range("A5:E5").Value = ""
$ws.Range("A6:E6").Value = ""

# --- Force the used range to extend through row 6 (A1:E6) without
#     disturbing formatting: a no-op border touch (clearing an already
#     absent border) registers every cell in A3:E6 in the sheet's used
#     range while keeping the default (style 0) formatting intact. ---
$ws.Range("A3:E6").Borders.LineStyle = -4142
